$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 are rewritten in place with the recomputed (new-TPM) values for the
# surviving Sending/Target cluster combinations (the "ECs" target cluster is
# dropped entirely). Rows 6-7 are removed below.

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.43424333333333
$ws.Range("H2").Value = 31.30273
$ws.Range("I2").Value = 0.9711091978791583
$ws.Range("J2").Value = 0.9711091978791584
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.212806333333333
$ws.Range("N2").Value = 24.638419
$ws.Range("O2").Value = 0.886321983523066
$ws.Range("P2").Value = 0.886321983523066
$ws.Range("Q2").Value = 85.6944197315411
$ws.Range("R2").Value = 771.2497775838699
$ws.Range("S2").Value = 0.8607154304817491
$ws.Range("T2").Value = 0.8607154304817493

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.43424333333333
$ws.Range("H3").Value = 31.30273
$ws.Range("I3").Value = 0.9711091978791583
$ws.Range("J3").Value = 0.9711091978791584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.053359333333333
$ws.Range("N3").Value = 3.160078
$ws.Range("O3").Value = 0.113678016476934
$ws.Range("P3").Value = 0.113678016476934
$ws.Range("Q3").Value = 10.99100760143778
$ws.Range("R3").Value = 98.91906841293999
$ws.Range("S3").Value = 0.1103937673974091
$ws.Range("T3").Value = 0.1103937673974091

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.310422
$ws.Range("H4").Value = 0.9312659999999999
$ws.Range("I4").Value = 0.02889080212084161
$ws.Range("J4").Value = 0.02889080212084161
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.212806333333333
$ws.Range("N4").Value = 24.638419
$ws.Range("O4").Value = 0.886321983523066
$ws.Range("P4").Value = 0.886321983523066
$ws.Range("Q4").Value = 2.549435767606
$ws.Range("R4").Value = 22.944921908454
$ws.Range("S4").Value = 0.02560655304131674
$ws.Range("T4").Value = 0.02560655304131674

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.310422
$ws.Range("H5").Value = 0.9312659999999999
$ws.Range("I5").Value = 0.02889080212084161
$ws.Range("J5").Value = 0.02889080212084161
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.053359333333333
$ws.Range("N5").Value = 3.160078
$ws.Range("O5").Value = 0.113678016476934
$ws.Range("P5").Value = 0.113678016476934
$ws.Range("Q5").Value = 0.326985910972
$ws.Range("R5").Value = 2.942873198748
$ws.Range("S5").Value = 0.003284249079524872
$ws.Range("T5").Value = 0.003284249079524872

# Drop the now-obsolete trailing rows (previously the MuSCs -> FAPs / MuSCs rows,
# which have been folded into rows 4-5 above); this also shrinks the used range
# to A1:T5 and lets the writer prune the now-unreferenced "ECs" shared string.
$ws.Rows("6:7").Delete()
